$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.706.59"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "3.171.66"
$ws.Range("E3").Value = "  -4.71%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "571.47"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "171.75"
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("D9").Value = "3.169.95"
$ws.Range("E9").Value = "  -4.65%  "
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").Value = "6.61"
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("E12").Value = "  -2.99%  "
$ws.Range("D13").Value = "3.723.17"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "27.22"
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("D16").Value = "65.642.32"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "3.176.76"
$ws.Range("E18").Value = "  -5.10%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "12.91"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("D21").Value = "362.01"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "7.27"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "68.82"
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("E25").Value = "  -3.54%  "
$ws.Range("D26").Value = "3.307.74"
$ws.Range("E26").Value = "  -4.83%  "
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  -5.18%  "
$ws.Range("D28").Value = "9.89"
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "5.39"
$ws.Range("E33").Value = "  -3.00%  "
$ws.Range("D34").Value = "22.06"
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("D35").Value = "6.62"
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "161.49"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "1.80"
$ws.Range("E40").Value = "  +3.85%  "
$ws.Range("D41").Value = "26.38"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").Value = "2.52"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").Value = "2.652.87"
$ws.Range("E43").Value = "  -1.63%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "39.77"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "327.94"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").Value = "23.77"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").Value = "  -0.35%  "
